# Weekly update: a new daily price record is inserted at the top of the
# Cilantro series (row 147), every existing record from row 147 downward
# is pushed one row down, and the record that used to be last (old row
# 235) becomes the new row 236.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D147").Value = 44606
$ws.Range("J147").Value = 60
$ws.Range("K147").Value = 7000
$ws.Range("L147").Value = 7000
$ws.Range("M147").Value = 7000
$ws.Range("P147").Value = 3500
$ws.Range("D148").Value = 44238
$ws.Range("D149").Value = 44242
$ws.Range("J149").Value = 40
$ws.Range("K149").Value = 4000
$ws.Range("L149").Value = 4000
$ws.Range("M149").Value = 4000
$ws.Range("P149").Value = 2000
$ws.Range("D150").Value = 44166
$ws.Range("J150").Value = 180
$ws.Range("K150").Value = 5000
$ws.Range("L150").Value = 5000
$ws.Range("M150").Value = 5000
$ws.Range("N150").Value = "`$/docena de atados (2 kilos)"
$ws.Range("O150").Value = "Región de La Araucanía"
$ws.Range("P150").Value = 2500
$ws.Range("Q150").Value = 2
$ws.Range("D151").Value = 44348
$ws.Range("J151").Value = 300
$ws.Range("K151").Value = 10000
$ws.Range("L151").Value = 10000
$ws.Range("M151").Value = 10000
$ws.Range("N151").Value = "`$/caja 36 atados"
$ws.Range("O151").Value = "Región Metropolitana"
$ws.Range("P151").Value = 278
$ws.Range("Q151").Value = 36
$ws.Range("D152").Value = 44223
$ws.Range("J152").Value = 20
$ws.Range("K152").Value = 4000
$ws.Range("L152").Value = 4000
$ws.Range("M152").Value = 4000
$ws.Range("P152").Value = 2000
$ws.Range("D153").Value = 44201
$ws.Range("K153").Value = 5000
$ws.Range("L153").Value = 5000
$ws.Range("M153").Value = 5000
$ws.Range("P153").Value = 2500
$ws.Range("D154").Value = 44579
$ws.Range("J154").Value = 80
$ws.Range("N154").Value = "`$/docena de atados (2 kilos)"
$ws.Range("O154").Value = "Región de La Araucanía"
$ws.Range("P154").Value = 5000
$ws.Range("Q154").Value = 2
$ws.Range("D155").Value = 44515
$ws.Range("K155").Value = 10000
$ws.Range("L155").Value = 10000
$ws.Range("M155").Value = 10000
$ws.Range("P155").Value = 278
$ws.Range("D156").Value = 44438
$ws.Range("J156").Value = 150
$ws.Range("K156").Value = 15000
$ws.Range("L156").Value = 15000
$ws.Range("M156").Value = 15000
$ws.Range("N156").Value = "`$/caja 36 atados"
$ws.Range("O156").Value = "Región Metropolitana"
$ws.Range("P156").Value = 417
$ws.Range("Q156").Value = 36
$ws.Range("D157").Value = 44249
$ws.Range("J157").Value = 40
$ws.Range("K157").Value = 6000
$ws.Range("L157").Value = 6000
$ws.Range("M157").Value = 6000
$ws.Range("P157").Value = 3000
$ws.Range("D158").Value = 44566
$ws.Range("J158").Value = 20
$ws.Range("K158").Value = 12000
$ws.Range("L158").Value = 12000
$ws.Range("M158").Value = 12000
$ws.Range("N158").Value = "`$/docena de atados (2 kilos)"
$ws.Range("O158").Value = "Región de La Araucanía"
$ws.Range("P158").Value = 6000
$ws.Range("Q158").Value = 2
$ws.Range("D159").Value = 44351
$ws.Range("J159").Value = 260
$ws.Range("D160").Value = 44494
$ws.Range("J160").Value = 100
$ws.Range("K160").Value = 10000
$ws.Range("L160").Value = 10000
$ws.Range("M160").Value = 10000
$ws.Range("N160").Value = "`$/caja 36 atados"
$ws.Range("O160").Value = "Región Metropolitana"
$ws.Range("P160").Value = 278
$ws.Range("Q160").Value = 36
$ws.Range("D161").Value = 44237
$ws.Range("K161").Value = 5000
$ws.Range("L161").Value = 5000
$ws.Range("M161").Value = 5000
$ws.Range("P161").Value = 2500
$ws.Range("D162").Value = 44300
$ws.Range("J162").Value = 20
$ws.Range("K162").Value = 6000
$ws.Range("L162").Value = 6000
$ws.Range("M162").Value = 6000
$ws.Range("P162").Value = 3000
$ws.Range("D163").Value = 44578
$ws.Range("J163").Value = 40
$ws.Range("K163").Value = 10000
$ws.Range("L163").Value = 10000
$ws.Range("M163").Value = 10000
$ws.Range("P163").Value = 5000
$ws.Range("D164").Value = 44225
$ws.Range("J164").Value = 100
$ws.Range("K164").Value = 4000
$ws.Range("L164").Value = 4000
$ws.Range("M164").Value = 4000
$ws.Range("N164").Value = "`$/docena de atados (2 kilos)"
$ws.Range("O164").Value = "Región de La Araucanía"
$ws.Range("P164").Value = 2000
$ws.Range("Q164").Value = 2
$ws.Range("D165").Value = 44411
$ws.Range("J165").Value = 300
$ws.Range("K165").Value = 13000
$ws.Range("L165").Value = 14000
$ws.Range("M165").Value = 13500
$ws.Range("N165").Value = "`$/caja 36 atados"
$ws.Range("O165").Value = "Región Metropolitana"
$ws.Range("P165").Value = 375
$ws.Range("Q165").Value = 36
$ws.Range("D166").Value = 44281
$ws.Range("J166").Value = 6000
$ws.Range("K166").Value = 20
$ws.Range("L166").Value = 20
$ws.Range("M166").Value = 20
$ws.Range("P166").Value = 10
$ws.Range("D167").Value = 44252
$ws.Range("J167").Value = 60
$ws.Range("D168").Value = 44271
$ws.Range("J168").Value = 100
$ws.Range("K168").Value = 7000
$ws.Range("L168").Value = 7000
$ws.Range("M168").Value = 7000
$ws.Range("P168").Value = 3500
$ws.Range("D169").Value = 44162
$ws.Range("J169").Value = 180
$ws.Range("K169").Value = 4000
$ws.Range("L169").Value = 5000
$ws.Range("M169").Value = 4500
$ws.Range("N169").Value = "`$/docena de atados (2 kilos)"
$ws.Range("O169").Value = "Región de La Araucanía"
$ws.Range("P169").Value = 2250
$ws.Range("Q169").Value = 2
$ws.Range("D170").Value = 44343
$ws.Range("J170").Value = 80
$ws.Range("D171").Value = 44516
$ws.Range("J171").Value = 250
$ws.Range("K171").Value = 10000
$ws.Range("L171").Value = 10000
$ws.Range("M171").Value = 10000
$ws.Range("N171").Value = "`$/caja 36 atados"
$ws.Range("O171").Value = "Región Metropolitana"
$ws.Range("P171").Value = 278
$ws.Range("Q171").Value = 36
$ws.Range("D172").Value = 44568
$ws.Range("J172").Value = 80
$ws.Range("K172").Value = 12000
$ws.Range("L172").Value = 12000
$ws.Range("M172").Value = 12000
$ws.Range("N172").Value = "`$/docena de atados (2 kilos)"
$ws.Range("O172").Value = "Región de La Araucanía"
$ws.Range("P172").Value = 6000
$ws.Range("Q172").Value = 2
$ws.Range("D173").Value = 44511
$ws.Range("J173").Value = 100
$ws.Range("K173").Value = 10000
$ws.Range("L173").Value = 10000
$ws.Range("M173").Value = 10000
$ws.Range("P173").Value = 278
$ws.Range("D174").Value = 44336
$ws.Range("J174").Value = 60
$ws.Range("K174").Value = 9000
$ws.Range("L174").Value = 9500
$ws.Range("M174").Value = 9250
$ws.Range("N174").Value = "`$/caja 36 atados"
$ws.Range("O174").Value = "Región Metropolitana"
$ws.Range("P174").Value = 257
$ws.Range("Q174").Value = 36
$ws.Range("D175").Value = 44231
$ws.Range("J175").Value = 20
$ws.Range("K175").Value = 4000
$ws.Range("L175").Value = 4000
$ws.Range("M175").Value = 4000
$ws.Range("N175").Value = "`$/docena de atados (2 kilos)"
$ws.Range("O175").Value = "Región de La Araucanía"
$ws.Range("P175").Value = 2000
$ws.Range("Q175").Value = 2
$ws.Range("D176").Value = 44400
$ws.Range("J176").Value = 280
$ws.Range("K176").Value = 13500
$ws.Range("L176").Value = 13500
$ws.Range("M176").Value = 13500
$ws.Range("P176").Value = 375
$ws.Range("D177").Value = 44334
$ws.Range("J177").Value = 300
$ws.Range("K177").Value = 9500
$ws.Range("M177").Value = 9750
$ws.Range("N177").Value = "`$/caja 36 atados"
$ws.Range("O177").Value = "Región Metropolitana"
$ws.Range("P177").Value = 271
$ws.Range("Q177").Value = 36
$ws.Range("D178").Value = 44573
$ws.Range("J178").Value = 20
$ws.Range("N178").Value = "`$/docena de atados (2 kilos)"
$ws.Range("O178").Value = "Región de La Araucanía"
$ws.Range("P178").Value = 5000
$ws.Range("Q178").Value = 2
$ws.Range("D179").Value = 44319
$ws.Range("K179").Value = 10000
$ws.Range("L179").Value = 10000
$ws.Range("M179").Value = 10000
$ws.Range("N179").Value = "`$/caja 36 atados"
$ws.Range("O179").Value = "Región Metropolitana"
$ws.Range("P179").Value = 278
$ws.Range("Q179").Value = 36
$ws.Range("D180").Value = 44280
$ws.Range("J180").Value = 60
$ws.Range("K180").Value = 6000
$ws.Range("L180").Value = 6000
$ws.Range("M180").Value = 6000
$ws.Range("N180").Value = "`$/docena de atados (2 kilos)"
$ws.Range("O180").Value = "Región de La Araucanía"
$ws.Range("P180").Value = 3000
$ws.Range("Q180").Value = 2
$ws.Range("D181").Value = 44362
$ws.Range("K181").Value = 9500
$ws.Range("L181").Value = 9500
$ws.Range("M181").Value = 9500
$ws.Range("P181").Value = 264
$ws.Range("D182").Value = 44365
$ws.Range("J182").Value = 300
$ws.Range("K182").Value = 10000
$ws.Range("L182").Value = 10000
$ws.Range("M182").Value = 10000
$ws.Range("N182").Value = "`$/caja 36 atados"
$ws.Range("O182").Value = "Región Metropolitana"
$ws.Range("P182").Value = 278
$ws.Range("Q182").Value = 36
$ws.Range("D183").Value = 44567
$ws.Range("J183").Value = 30
$ws.Range("N183").Value = "`$/docena de atados (2 kilos)"
$ws.Range("O183").Value = "Región de La Araucanía"
$ws.Range("P183").Value = 6000
$ws.Range("Q183").Value = 2
$ws.Range("D184").Value = 44473
$ws.Range("J184").Value = 200
$ws.Range("K184").Value = 12000
$ws.Range("L184").Value = 12000
$ws.Range("M184").Value = 12000
$ws.Range("P184").Value = 333
$ws.Range("D185").Value = 44357
$ws.Range("J185").Value = 100
$ws.Range("K185").Value = 9500
$ws.Range("L185").Value = 9500
$ws.Range("M185").Value = 9500
$ws.Range("N185").Value = "`$/caja 36 atados"
$ws.Range("O185").Value = "Región Metropolitana"
$ws.Range("P185").Value = 264
$ws.Range("Q185").Value = 36
$ws.Range("D186").Value = 44581
$ws.Range("J186").Value = 50
$ws.Range("K186").Value = 12000
$ws.Range("L186").Value = 12000
$ws.Range("M186").Value = 12000
$ws.Range("P186").Value = 6000
$ws.Range("D187").Value = 44537
$ws.Range("J187").Value = 160
$ws.Range("K187").Value = 6000
$ws.Range("L187").Value = 6000
$ws.Range("M187").Value = 6000
$ws.Range("P187").Value = 3000
$ws.Range("D188").Value = 44553
$ws.Range("J188").Value = 80
$ws.Range("K188").Value = 7000
$ws.Range("L188").Value = 7000
$ws.Range("M188").Value = 7000
$ws.Range("N188").Value = "`$/docena de atados (2 kilos)"
$ws.Range("O188").Value = "Región de La Araucanía"
$ws.Range("P188").Value = 3500
$ws.Range("Q188").Value = 2
$ws.Range("D189").Value = 44490
$ws.Range("J189").Value = 100
$ws.Range("K189").Value = 10000
$ws.Range("L189").Value = 10000
$ws.Range("M189").Value = 10000
$ws.Range("P189").Value = 278
$ws.Range("D190").Value = 44397
$ws.Range("J190").Value = 180
$ws.Range("K190").Value = 14500
$ws.Range("L190").Value = 14500
$ws.Range("M190").Value = 14500
$ws.Range("P190").Value = 403
$ws.Range("D191").Value = 44462
$ws.Range("J191").Value = 150
$ws.Range("L191").Value = 12000
$ws.Range("M191").Value = 12000
$ws.Range("P191").Value = 333
$ws.Range("D192").Value = 44446
$ws.Range("K192").Value = 12000
$ws.Range("L192").Value = 13000
$ws.Range("M192").Value = 12500
$ws.Range("P192").Value = 347
$ws.Range("D193").Value = 44421
$ws.Range("J193").Value = 300
$ws.Range("K193").Value = 13000
$ws.Range("L193").Value = 15000
$ws.Range("M193").Value = 14000
$ws.Range("P193").Value = 389
$ws.Range("D194").Value = 44329
$ws.Range("K194").Value = 11000
$ws.Range("L194").Value = 11000
$ws.Range("M194").Value = 11000
$ws.Range("N194").Value = "`$/caja 36 atados"
$ws.Range("O194").Value = "Región Metropolitana"
$ws.Range("P194").Value = 306
$ws.Range("Q194").Value = 36
$ws.Range("D195").Value = 44208
$ws.Range("J195").Value = 100
$ws.Range("K195").Value = 5000
$ws.Range("L195").Value = 5000
$ws.Range("M195").Value = 5000
$ws.Range("N195").Value = "`$/docena de atados (2 kilos)"
$ws.Range("O195").Value = "Región de La Araucanía"
$ws.Range("P195").Value = 2500
$ws.Range("Q195").Value = 2
$ws.Range("D196").Value = 44355
$ws.Range("J196").Value = 290
$ws.Range("K196").Value = 9500
$ws.Range("L196").Value = 9500
$ws.Range("M196").Value = 9500
$ws.Range("N196").Value = "`$/caja 36 atados"
$ws.Range("O196").Value = "Región Metropolitana"
$ws.Range("P196").Value = 264
$ws.Range("Q196").Value = 36
$ws.Range("D197").Value = 44530
$ws.Range("J197").Value = 200
$ws.Range("K197").Value = 4000
$ws.Range("L197").Value = 5000
$ws.Range("M197").Value = 4500
$ws.Range("N197").Value = "`$/docena de atados (2 kilos)"
$ws.Range("O197").Value = "Región de La Araucanía"
$ws.Range("P197").Value = 2250
$ws.Range("Q197").Value = 2
$ws.Range("D198").Value = 44483
$ws.Range("J198").Value = 150
$ws.Range("K198").Value = 10000
$ws.Range("L198").Value = 10000
$ws.Range("M198").Value = 10000
$ws.Range("P198").Value = 278
$ws.Range("D199").Value = 44294
$ws.Range("J199").Value = 120
$ws.Range("K199").Value = 12000
$ws.Range("L199").Value = 12000
$ws.Range("M199").Value = 12000
$ws.Range("N199").Value = "`$/caja 36 atados"
$ws.Range("O199").Value = "Región Metropolitana"
$ws.Range("P199").Value = 333
$ws.Range("Q199").Value = 36
$ws.Range("D200").Value = 44264
$ws.Range("J200").Value = 80
$ws.Range("K200").Value = 7000
$ws.Range("L200").Value = 7000
$ws.Range("M200").Value = 7000
$ws.Range("N200").Value = "`$/docena de atados (2 kilos)"
$ws.Range("O200").Value = "Región de La Araucanía"
$ws.Range("P200").Value = 3500
$ws.Range("Q200").Value = 2
$ws.Range("D201").Value = 44396
$ws.Range("J201").Value = 70
$ws.Range("K201").Value = 12000
$ws.Range("L201").Value = 12000
$ws.Range("M201").Value = 12000
$ws.Range("N201").Value = "`$/caja 36 atados"
$ws.Range("O201").Value = "Región Metropolitana"
$ws.Range("P201").Value = 333
$ws.Range("Q201").Value = 36
$ws.Range("D202").Value = 44232
$ws.Range("J202").Value = 100
$ws.Range("K202").Value = 4000
$ws.Range("L202").Value = 4000
$ws.Range("M202").Value = 4000
$ws.Range("P202").Value = 2000
$ws.Range("D203").Value = 44279
$ws.Range("J203").Value = 30
$ws.Range("K203").Value = 6000
$ws.Range("L203").Value = 6000
$ws.Range("M203").Value = 6000
$ws.Range("N203").Value = "`$/docena de atados (2 kilos)"
$ws.Range("O203").Value = "Región de La Araucanía"
$ws.Range("P203").Value = 3000
$ws.Range("Q203").Value = 2
$ws.Range("D204").Value = 44330
$ws.Range("J204").Value = 300
$ws.Range("K204").Value = 10000
$ws.Range("L204").Value = 10000
$ws.Range("M204").Value = 10000
$ws.Range("P204").Value = 278
$ws.Range("D205").Value = 44504
$ws.Range("J205").Value = 150
$ws.Range("K205").Value = 9000
$ws.Range("L205").Value = 9000
$ws.Range("M205").Value = 9000
$ws.Range("N205").Value = "`$/caja 36 atados"
$ws.Range("O205").Value = "Región Metropolitana"
$ws.Range("P205").Value = 250
$ws.Range("Q205").Value = 36
$ws.Range("D206").Value = 44572
$ws.Range("J206").Value = 120
$ws.Range("K206").Value = 10000
$ws.Range("L206").Value = 10000
$ws.Range("M206").Value = 10000
$ws.Range("P206").Value = 5000
$ws.Range("D207").Value = 44257
$ws.Range("K207").Value = 8000
$ws.Range("L207").Value = 8000
$ws.Range("M207").Value = 8000
$ws.Range("N207").Value = "`$/docena de atados (2 kilos)"
$ws.Range("O207").Value = "Región de La Araucanía"
$ws.Range("P207").Value = 4000
$ws.Range("Q207").Value = 2
$ws.Range("D208").Value = 44301
$ws.Range("J208").Value = 100
$ws.Range("K208").Value = 13000
$ws.Range("L208").Value = 13000
$ws.Range("M208").Value = 13000
$ws.Range("P208").Value = 361
$ws.Range("D209").Value = 44370
$ws.Range("J209").Value = 50
$ws.Range("D210").Value = 44487
$ws.Range("J210").Value = 150
$ws.Range("K210").Value = 10000
$ws.Range("L210").Value = 10000
$ws.Range("M210").Value = 10000
$ws.Range("N210").Value = "`$/caja 36 atados"
$ws.Range("O210").Value = "Región Metropolitana"
$ws.Range("P210").Value = 278
$ws.Range("Q210").Value = 36
$ws.Range("D211").Value = 44174
$ws.Range("J211").Value = 25
$ws.Range("D212").Value = 44200
$ws.Range("J212").Value = 50
$ws.Range("K212").Value = 5000
$ws.Range("L212").Value = 5000
$ws.Range("M212").Value = 5000
$ws.Range("N212").Value = "`$/docena de atados (2 kilos)"
$ws.Range("O212").Value = "Región de La Araucanía"
$ws.Range("P212").Value = 2500
$ws.Range("Q212").Value = 2
$ws.Range("D213").Value = 44385
$ws.Range("K213").Value = 12000
$ws.Range("L213").Value = 12000
$ws.Range("M213").Value = 12000
$ws.Range("N213").Value = "`$/caja 36 atados"
$ws.Range("O213").Value = "Región Metropolitana"
$ws.Range("P213").Value = 333
$ws.Range("Q213").Value = 36
$ws.Range("D214").Value = 44236
$ws.Range("K214").Value = 5000
$ws.Range("L214").Value = 5000
$ws.Range("M214").Value = 5000
$ws.Range("P214").Value = 2500
$ws.Range("D215").Value = 44221
$ws.Range("J215").Value = 100
$ws.Range("K215").Value = 4000
$ws.Range("L215").Value = 4000
$ws.Range("M215").Value = 4000
$ws.Range("N215").Value = "`$/docena de atados (2 kilos)"
$ws.Range("O215").Value = "Región de La Araucanía"
$ws.Range("P215").Value = 2000
$ws.Range("Q215").Value = 2
$ws.Range("D216").Value = 44413
$ws.Range("J216").Value = 140
$ws.Range("K216").Value = 15000
$ws.Range("L216").Value = 16000
$ws.Range("M216").Value = 15500
$ws.Range("N216").Value = "`$/caja 36 atados"
$ws.Range("O216").Value = "Región Metropolitana"
$ws.Range("P216").Value = 431
$ws.Range("Q216").Value = 36
$ws.Range("D217").Value = 44272
$ws.Range("J217").Value = 20
$ws.Range("K217").Value = 7000
$ws.Range("L217").Value = 7000
$ws.Range("M217").Value = 7000
$ws.Range("P217").Value = 3500
$ws.Range("D218").Value = 44229
$ws.Range("J218").Value = 200
$ws.Range("K218").Value = 4000
$ws.Range("L218").Value = 5000
$ws.Range("M218").Value = 4500
$ws.Range("N218").Value = "`$/docena de atados (2 kilos)"
$ws.Range("O218").Value = "Región de La Araucanía"
$ws.Range("P218").Value = 2250
$ws.Range("Q218").Value = 2
$ws.Range("J219").Value = 180
$ws.Range("K219").Value = 13000
$ws.Range("L219").Value = 13500
$ws.Range("M219").Value = 13250
$ws.Range("N219").Value = "`$/caja 36 atados"
$ws.Range("O219").Value = "Región Metropolitana"
$ws.Range("P219").Value = 368
$ws.Range("Q219").Value = 36
$ws.Range("D220").Value = 44299
$ws.Range("J220").Value = 80
$ws.Range("D221").Value = 44312
$ws.Range("J221").Value = 20
$ws.Range("K221").Value = 6000
$ws.Range("L221").Value = 6000
$ws.Range("M221").Value = 6000
$ws.Range("N221").Value = "`$/docena de atados (2 kilos)"
$ws.Range("O221").Value = "Región de La Araucanía"
$ws.Range("P221").Value = 3000
$ws.Range("Q221").Value = 2
$ws.Range("D222").Value = 44399
$ws.Range("J222").Value = 110
$ws.Range("K222").Value = 12000
$ws.Range("L222").Value = 12000
$ws.Range("M222").Value = 12000
$ws.Range("N222").Value = "`$/caja 36 atados"
$ws.Range("O222").Value = "Región Metropolitana"
$ws.Range("P222").Value = 333
$ws.Range("Q222").Value = 36
$ws.Range("D223").Value = 44167
$ws.Range("J223").Value = 30
$ws.Range("K223").Value = 5000
$ws.Range("L223").Value = 5000
$ws.Range("M223").Value = 5000
$ws.Range("P223").Value = 2500
$ws.Range("D224").Value = 44277
$ws.Range("J224").Value = 90
$ws.Range("K224").Value = 6000
$ws.Range("L224").Value = 6000
$ws.Range("M224").Value = 6000
$ws.Range("P224").Value = 3000
$ws.Range("D225").Value = 44258
$ws.Range("J225").Value = 20
$ws.Range("K225").Value = 8000
$ws.Range("L225").Value = 8000
$ws.Range("M225").Value = 8000
$ws.Range("N225").Value = "`$/docena de atados (2 kilos)"
$ws.Range("O225").Value = "Región de La Araucanía"
$ws.Range("P225").Value = 4000
$ws.Range("Q225").Value = 2
$ws.Range("D226").Value = 44390
$ws.Range("J226").Value = 300
$ws.Range("K226").Value = 12000
$ws.Range("L226").Value = 12000
$ws.Range("M226").Value = 12000
$ws.Range("P226").Value = 333
$ws.Range("D227").Value = 44349
$ws.Range("J227").Value = 37
$ws.Range("K227").Value = 9500
$ws.Range("L227").Value = 9500
$ws.Range("M227").Value = 9500
$ws.Range("P227").Value = 264
$ws.Range("J228").Value = 200
$ws.Range("K228").Value = 14000
$ws.Range("L228").Value = 14000
$ws.Range("M228").Value = 14000
$ws.Range("N228").Value = "`$/caja 36 atados"
$ws.Range("O228").Value = "Región Metropolitana"
$ws.Range("P228").Value = 389
$ws.Range("Q228").Value = 36
$ws.Range("D229").Value = 44285
$ws.Range("J229").Value = 100
$ws.Range("K229").Value = 8000
$ws.Range("L229").Value = 8000
$ws.Range("M229").Value = 8000
$ws.Range("N229").Value = "`$/docena de atados (2 kilos)"
$ws.Range("O229").Value = "Región de La Araucanía"
$ws.Range("P229").Value = 4000
$ws.Range("Q229").Value = 2
$ws.Range("D230").Value = 44498
$ws.Range("J230").Value = 350
$ws.Range("K230").Value = 10000
$ws.Range("L230").Value = 10000
$ws.Range("M230").Value = 10000
$ws.Range("P230").Value = 278
$ws.Range("D231").Value = 44418
$ws.Range("J231").Value = 300
$ws.Range("K231").Value = 14000
$ws.Range("M231").Value = 14500
$ws.Range("P231").Value = 403
$ws.Range("J232").Value = 70
$ws.Range("K232").Value = 15000
$ws.Range("L232").Value = 15000
$ws.Range("M232").Value = 15000
$ws.Range("N232").Value = "`$/caja 36 atados"
$ws.Range("O232").Value = "Región Metropolitana"
$ws.Range("P232").Value = 417
$ws.Range("Q232").Value = 36
$ws.Range("D233").Value = 44595
$ws.Range("J233").Value = 40
$ws.Range("K233").Value = 11000
$ws.Range("L233").Value = 11000
$ws.Range("M233").Value = 11000
$ws.Range("N233").Value = "`$/docena de atados (2 kilos)"
$ws.Range("O233").Value = "Región de La Araucanía"
$ws.Range("P233").Value = 5500
$ws.Range("Q233").Value = 2
$ws.Range("D234").Value = 44335
$ws.Range("J234").Value = 30
$ws.Range("K234").Value = 9500
$ws.Range("L234").Value = 9500
$ws.Range("M234").Value = 9500
$ws.Range("N234").Value = "`$/caja 36 atados"
$ws.Range("O234").Value = "Región Metropolitana"
$ws.Range("P234").Value = 264
$ws.Range("Q234").Value = 36
$ws.Range("D235").Value = 44552
$ws.Range("J235").Value = 20
$ws.Range("K235").Value = 7000
$ws.Range("L235").Value = 7000
$ws.Range("M235").Value = 7000
$ws.Range("P235").Value = 3500
$ws.Range("A236").Value = 4
$ws.Range("B236").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C236").Value = "Los Lagos"
$ws.Range("D236").Value = 44544
$ws.Range("E236").Value = 10
$ws.Range("F236").Value = 100112040
$ws.Range("G236").Value = "Cilantro"
$ws.Range("H236").Value = "Sin especificar"
$ws.Range("I236").Value = "Primera"
$ws.Range("J236").Value = 180
$ws.Range("K236").Value = 6000
$ws.Range("L236").Value = 6000
$ws.Range("M236").Value = 6000
$ws.Range("N236").Value = "`$/docena de atados (2 kilos)"
$ws.Range("O236").Value = "Región de La Araucanía"
$ws.Range("P236").Value = 3000
$ws.Range("Q236").Value = 2
$ws.Range("R236").Value = "Hortaliza"

# Match date style for new row 236 (column D) to the rest of the date column
$ws.Range("D236").NumberFormat = $ws.Range("D235").NumberFormat

